# Fruta / hortaliza, semanal
# Inserts a new weekly price record as row 146 in the "Granada" sheet
# (Vega Modelo de Temuco), shifting the existing rows 146:183 down to
# 147:184 and growing the used range from A1:T183 to A1:T184.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 146 - this shifts rows 146:183 down to 147:184
# and extends the sheet dimension automatically.
$ws.Rows.Item(146).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A146").Value = 10
$ws.Range("B146").Value = "Vega Modelo de Temuco"
$ws.Range("C146").Value = "La Araucanía"
$ws.Range("D146").Value = 44841
$ws.Range("E146").Value = 9
$ws.Range("F146").Value = "Fruta"
$ws.Range("G146").Value = 100104
$ws.Range("H146").Value = "Frutos de pepita"
$ws.Range("I146").Value = 100104001
$ws.Range("J146").Value = "Granada"
$ws.Range("K146").Value = "Wonderfull"
$ws.Range("L146").Value = "Primera"
$ws.Range("M146").Value = 65
$ws.Range("N146").Value = 18000
$ws.Range("O146").Value = 18000
$ws.Range("P146").Value = 18000
$ws.Range("Q146").Value = "$/bandeja 15 kilos granel"
$ws.Range("R146").Value = "Provincia de Limarí"
$ws.Range("S146").Value = 1200
$ws.Range("T146").Value = 15
